# "added EAIAK 150 Mark" — fill in the previously-empty EAIAK value/STD
# pair for the Temperature = 150 row (row 8) on Sheet1, matching the
# pattern already used by the other data rows (D2:D9 / E2:E9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D8 = EAIAK value, E8 = EAIAK STD for the 150-degree mark.
$ws.Range("D8").Value = 0.0873083333333
$ws.Range("E8").Value = 0.0375345609258

# The other populated cells in these columns carry the sheet's default
# (unstyled) formatting rather than the workbook's generic "applied font"
# style that empty D8/E8 previously had. Re-asserting the default font
# nudges the new cells onto that same plain style, matching the rest of
# the column instead of inheriting the old empty-cell style.
$ws.Range("D8:E8").Font.Name = "Arial"

# Restore/update the active selection, as the author's session ended
# with E11 selected rather than the original E13.
$ws.Range("E11").Select() | Out-Null
